# Covid_Dominican / Dominican_Covid.xlsx
# Commit: "Latest date 5 February"
#
# Fills in the next two weekly report dates (29 Jan 2021 / serial 44225 and
# 5 Feb 2021 / serial 44232) on both the national summary sheet
# (Fallecido_Recuperado) and the per-province weekly sheet
# (Provincias_Semanal), and stubs in the following week's date
# (12 Feb 2021 / serial 44239) with no data yet - matching the existing
# pattern already present in the workbook for an as-yet-unreported week.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Fallecido_Recuperado": national cumulative totals, one row/week
# ---------------------------------------------------------------------
$wsNat = $wb.Worksheets.Item("Fallecido_Recuperado")

# Row 46 (29-Jan-2021 / 44225) already has its date; fill the totals.
$wsNat.Cells.Item(46, 2).Value = 212553
$wsNat.Cells.Item(46, 3).Value = 2646
$wsNat.Cells.Item(46, 4).Value = 155867

# Row 47 (5-Feb-2021 / 44232): new week, full data.
$wsNat.Cells.Item(47, 1).Value = 44232
$wsNat.Cells.Item(47, 2).Value = 222148
$wsNat.Cells.Item(47, 3).Value = 2801
$wsNat.Cells.Item(47, 4).Value = 165659

# Row 48 (12-Feb-2021 / 44239): date only, not yet reported.
$wsNat.Cells.Item(48, 1).Value = 44239

# New date cells need the same date number format as the existing column A.
$wsNat.Range("A46").Copy()
$wsNat.Range("A47:A48").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet "Provincias_Semanal": per-province weekly data, 32 rows/week
# ---------------------------------------------------------------------
$wsProv = $wb.Worksheets.Item("Provincias_Semanal")

$provData = @(
        @{ Row=1410; A=44225; B="Distrito Nacional"; C=5727.07; D=443 },
        @{ Row=1411; A=44225; B="Azua"; C=1113.13; D=29 },
        @{ Row=1412; A=44225; B="Baoruco"; C=1228.65; D=10 },
        @{ Row=1413; A=44225; B="Barahona"; C=1282.4; D=23 },
        @{ Row=1414; A=44225; B="Dajabon"; C=1223.16; D=8 },
        @{ Row=1415; A=44225; B="Duarte"; C=1958.76; D=151 },
        @{ Row=1416; A=44225; B="Elias Pina"; C=567.81; D=5 },
        @{ Row=1417; A=44225; B="El Seibo"; C=811.7; D=7 },
        @{ Row=1418; A=44225; B="Espaillat"; C=1870.37; D=100 },
        @{ Row=1419; A=44225; B="Independencia"; C=1489.11; D=9 },
        @{ Row=1420; A=44225; B="La Altagracia"; C=2255.76; D=52 },
        @{ Row=1421; A=44225; B="La Romana"; C=2263.05; D=119 },
        @{ Row=1422; A=44225; B="La Vega"; C=2156.87; D=135 },
        @{ Row=1423; A=44225; B="Maria Trinidad Sanchez"; C=1924.92; D=18 },
        @{ Row=1424; A=44225; B="Monte Cristi"; C=789.11; D=15 },
        @{ Row=1425; A=44225; B="Pedernales"; C=1723; D=3 },
        @{ Row=1426; A=44225; B="Peravia"; C=892.45; D=47 },
        @{ Row=1427; A=44225; B="Puerto Plata"; C=2084.2; D=137 },
        @{ Row=1428; A=44225; B="Hermanas Mirabal"; C=2039.11; D=28 },
        @{ Row=1429; A=44225; B="Samana"; C=678.5; D=4 },
        @{ Row=1430; A=44225; B="San Cristobal"; C=1005.45; D=127 },
        @{ Row=1431; A=44225; B="San Juan"; C=1407.98; D=47 },
        @{ Row=1432; A=44225; B="San Pedro de Macoris"; C=879.77; D=41 },
        @{ Row=1433; A=44225; B="Sanchez Ramirez"; C=1855.56; D=25 },
        @{ Row=1434; A=44225; B="Santiago"; C=2255.52; D=407 },
        @{ Row=1435; A=44225; B="Santiago Rodriguez"; C=1650.19; D=10 },
        @{ Row=1436; A=44225; B="Valverde"; C=1126.08; D=31 },
        @{ Row=1437; A=44225; B="Monsenor Nouel"; C=1829.26; D=39 },
        @{ Row=1438; A=44225; B="Monte Plata"; C=533.94; D=26 },
        @{ Row=1439; A=44225; B="Hato Mayor"; C=678.62; D=13 },
        @{ Row=1440; A=44225; B="San Jose de Ocoa"; C=1373.04; D=15 },
        @{ Row=1441; A=44225; B="Santo Domingo"; C=1426.05; D=522 },
        @{ Row=1442; A=44232; B="Distrito Nacional"; C=6006.22; D=448 },
        @{ Row=1443; A=44232; B="Azua"; C=1146.88; D=29 },
        @{ Row=1444; A=44232; B="Baoruco"; C=1262.31; D=10 },
        @{ Row=1445; A=44232; B="Barahona"; C=1318.34; D=27 },
        @{ Row=1446; A=44232; B="Dajabon"; C=1314.94; D=10 },
        @{ Row=1447; A=44232; B="Duarte"; C=2009.61; D=165 },
        @{ Row=1448; A=44232; B="Elias Pina"; C=607.24; D=5 },
        @{ Row=1449; A=44232; B="El Seibo"; C=853.4; D=9 },
        @{ Row=1450; A=44232; B="Espaillat"; C=1962.1; D=104 },
        @{ Row=1451; A=44232; B="Independencia"; C=1533.62; D=10 },
        @{ Row=1452; A=44232; B="La Altagracia"; C=2355.93; D=56 },
        @{ Row=1453; A=44232; B="La Romana"; C=2365.03; D=120 },
        @{ Row=1454; A=44232; B="La Vega"; C=2245.13; D=152 },
        @{ Row=1455; A=44232; B="Maria Trinidad Sanchez"; C=2012.8; D=21 },
        @{ Row=1456; A=44232; B="Monte Cristi"; C=830.91; D=18 },
        @{ Row=1457; A=44232; B="Pedernales"; C=1725.86; D=3 },
        @{ Row=1458; A=44232; B="Peravia"; C=913.72; D=47 },
        @{ Row=1459; A=44232; B="Puerto Plata"; C=2154.13; D=143 },
        @{ Row=1460; A=44232; B="Hermanas Mirabal"; C=2199.72; D=31 },
        @{ Row=1461; A=44232; B="Samana"; C=718.62; D=4 },
        @{ Row=1462; A=44232; B="San Cristobal"; C=1035.41; D=134 },
        @{ Row=1463; A=44232; B="San Juan"; C=1463.45; D=48 },
        @{ Row=1464; A=44232; B="San Pedro de Macoris"; C=930.6; D=42 },
        @{ Row=1465; A=44232; B="Sanchez Ramirez"; C=1906.85; D=31 },
        @{ Row=1466; A=44232; B="Santiago"; C=2406.88; D=443 },
        @{ Row=1467; A=44232; B="Santiago Rodriguez"; C=1772.43; D=10 },
        @{ Row=1468; A=44232; B="Valverde"; C=1197.94; D=35 },
        @{ Row=1469; A=44232; B="Monsenor Nouel"; C=1879.76; D=42 },
        @{ Row=1470; A=44232; B="Monte Plata"; C=554.88; D=27 },
        @{ Row=1471; A=44232; B="Hato Mayor"; C=705.44; D=13 },
        @{ Row=1472; A=44232; B="San Jose de Ocoa"; C=1404.28; D=15 },
        @{ Row=1473; A=44232; B="Santo Domingo"; C=1489.88; D=549 }
)

foreach ($item in $provData) {
    $wsProv.Cells.Item($item.Row, 1).Value = $item.A
    $wsProv.Cells.Item($item.Row, 2).Value = $item.B
    $wsProv.Cells.Item($item.Row, 3).Value = $item.C
    $wsProv.Cells.Item($item.Row, 4).Value = $item.D
}

# Rows 1411-1473 are brand new rows; give column A the same date format
# already used on A1410 (and throughout the rest of the column).
$wsProv.Range("A1410").Copy()
$wsProv.Range("A1411:A1473").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Selections / window state: mirror the author's final cursor position
# on each sheet, while leaving Provincias_Semanal as the active tab.
# ---------------------------------------------------------------------
$wsNat.Range("B48").Select()
$wsProv.Range("A1474").Select()
